# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1) Metadata sheet: bump the "Date" property value.
# 2) Elements sheet: add a new "Mapping: Spécification métier vers
#    l'extension ROR CommuneCog" column (AL) and fill in the mapping
#    value ("communeCog") for the Extension.value[x] row.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 : refresh the publication Date value -------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2) Elements sheet : new mapping column ---------------------------------
$ws = $wb.Worksheets.Item("Elements")

$lastCol = 37   # AK = existing "Mapping: RIM Mapping" column
$newCol = 38    # AL = new mapping column
$lastRow = 6    # rows 2..6 hold data (row 1 is the header)

# Header cell: copy formatting from the previous mapping header (AK1) so the
# new column matches the bold/filled header style, then set its caption.
$ws.Cells.Item(1, $lastCol).Copy()
$ws.Cells.Item(1, $newCol).PasteSpecial(-4122)
$ws.Cells.Item(1, $newCol).Value = "Mapping: Spécification métier vers l'extension ROR CommuneCog"

# Data cells: copy formatting from the existing mapping column for every data
# row so the new column's blank cells keep the same borders/alignment.
$ws.Range($ws.Cells.Item(2, $lastCol), $ws.Cells.Item($lastRow, $lastCol)).Copy()
$ws.Range($ws.Cells.Item(2, $newCol), $ws.Cells.Item($lastRow, $newCol)).PasteSpecial(-4122)

# Only the Extension.value[x] row (row 6) gets an actual mapping value.
$ws.Cells.Item(6, $newCol).Value = "communeCog"

# Size the new column similarly to the other wide, best-fit text columns.
$ws.Columns.Item($newCol).ColumnWidth = 68.3
